# Update weekly Haba price data rows (3-13) for Mapocho Venta Directa de Santiago
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44446
$ws.Range("J3").Value = 15
$ws.Range("K3").Value = 13000
$ws.Range("L3").Value = 13000
$ws.Range("M3").Value = 13000
$ws.Range("O3").Value = "Provincia de Limarí"
$ws.Range("P3").Value = 520

# Row 4
$ws.Range("D4").Value = 44467
$ws.Range("J4").Value = 35
$ws.Range("K4").Value = 12000
$ws.Range("L4").Value = 12000
$ws.Range("M4").Value = 12000
$ws.Range("O4").Value = "Provincia de Limarí"
$ws.Range("P4").Value = 480

# Row 6
$ws.Range("D6").Value = 44453
$ws.Range("J6").Value = 55
$ws.Range("K6").Value = 14000
$ws.Range("L6").Value = 15000
$ws.Range("M6").Value = 14455
$ws.Range("P6").Value = 578

# Row 7
$ws.Range("D7").Value = 44432
$ws.Range("J7").Value = 15
$ws.Range("K7").Value = 14000
$ws.Range("L7").Value = 14000
$ws.Range("M7").Value = 14000
$ws.Range("O7").Value = "Provincia del Elquí"
$ws.Range("P7").Value = 560

# Row 8
$ws.Range("D8").Value = 44340
$ws.Range("J8").Value = 25
$ws.Range("K8").Value = 15000
$ws.Range("L8").Value = 15000
$ws.Range("M8").Value = 15000
$ws.Range("P8").Value = 600

# Row 9
$ws.Range("D9").Value = 44449
$ws.Range("J9").Value = 30
$ws.Range("K9").Value = 16000
$ws.Range("L9").Value = 16000
$ws.Range("M9").Value = 16000
$ws.Range("P9").Value = 640

# Row 10
$ws.Range("O10").Value = "Provincia de Limarí"

# Row 11
$ws.Range("D11").Value = 44435
$ws.Range("J11").Value = 15
$ws.Range("K11").Value = 14000
$ws.Range("L11").Value = 14000
$ws.Range("M11").Value = 14000
$ws.Range("O11").Value = "Provincia del Elquí"
$ws.Range("P11").Value = 560

# Row 12
$ws.Range("D12").Value = 44418
$ws.Range("J12").Value = 12
$ws.Range("K12").Value = 15000
$ws.Range("L12").Value = 15000
$ws.Range("M12").Value = 15000
$ws.Range("P12").Value = 600

# Row 13
$ws.Range("D13").Value = 44421
$ws.Range("J13").Value = 20
$ws.Range("K13").Value = 15000
$ws.Range("L13").Value = 15000
$ws.Range("M13").Value = 15000
$ws.Range("P13").Value = 600
